# edit.ps1
# Converts the docxtemplater-style "{{ }}" / "{% ... %}" Jinja-ish
# placeholders used by this A4 label template into single-brace
# "{ }" placeholders, and drops the (hidden) "{% for %}" / "{% endfor %}"
# loop-control runs that docxtemplater no longer needs for this engine.
#
# Summary of the edit:
#   1. The hidden paragraph that holds "{% for e in tbl_contents %} "
#      loses its two runs (the paragraph itself - and its formatting -
#      stays in place, now empty).
#   2. NAME cell: "{{e_name}}{{e_surname}}" -> "{e_name} {e_surname}"
#      (double braces collapse to single braces and the two fields
#      are separated by a literal space instead of being glued together).
#   3. ADDRESS / Email / Phone cells: "{{e_xxx}}" -> "{e_xxx}".
#   4. The hidden paragraph that holds "{% endfor %}" loses its three
#      runs (the paragraph itself stays in place, now empty).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: replace the first remaining occurrence of $findText in the
# whole document with $replText. We deliberately use Replace:=1
# (wdReplaceOne, not wdReplaceAll) and Wrap:=0 (wdFindStop) and always
# search $d.Content (the whole story) because this engine's Find
# always scans from the top of the story - it does not honour a
# sub-range's Start as a search-from position. Using ReplaceOne plus
# a tight, unambiguous search string (that lives fully inside a single
# run) lets us target one specific occurrence deterministically and
# keeps every other run's formatting untouched.
# ---------------------------------------------------------------------
function Replace-FirstOccurrence {
    param(
        [string]$findText,
        [string]$replText
    )
    $ok = $d.Content.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 0, $false, $replText, 1)
    if (-not $ok) {
        Write-Host "WARNING: could not find '$findText'"
    }
    return $ok
}

# ---------------------------------------------------------------------
# Helper: locate the hidden (vanish) hand-authored "{% ... %}" control
# paragraphs by scanning every paragraph's underlying OOXML for a
# marker substring, then wipe their runs via InsertXML while leaving
# the paragraph mark / pPr (and therefore the paragraph's formatting)
# completely intact. Range.Text/Range.Delete do not address hidden
# runs individually in this engine, so InsertXML (which replaces the
# full contents of the exact Range it's called on) is the reliable way
# to drop runs without disturbing the paragraph mark.
# ---------------------------------------------------------------------
function Clear-HiddenControlParagraph {
    param(
        [string]$marker
    )
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Font.Hidden) {
            $xml = $p.Range.WordOpenXML
            if ($xml -like "*$marker*") {
                $pPrStart = $xml.IndexOf("<w:pPr")
                if ($pPrStart -ge 0) {
                    $pPrEnd = $xml.IndexOf("</w:pPr>", $pPrStart) + "</w:pPr>".Length
                    $pPr = $xml.Substring($pPrStart, $pPrEnd - $pPrStart)
                } else {
                    $pPr = ""
                }
                $newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $pPr + '</w:p>'
                $p.Range.InsertXML($newXml)
                return $true
            }
        }
    }
    Write-Host "WARNING: could not find hidden paragraph containing '$marker'"
    return $false
}

# -----------------------------------------------------------------
# 1. Drop the "{% for e in tbl_contents %} " hidden runs.
# -----------------------------------------------------------------
Clear-HiddenControlParagraph "tbl_contents"

# -----------------------------------------------------------------
# 2. NAME cell: split "{{e_name}}{{e_surname}}" into
#    "{e_name} {e_surname}" (note the new run boundary + space).
#    This reshapes run boundaries (a brand-new run is introduced),
#    so we replace the whole paragraph's XML rather than doing a
#    plain text Find/Replace.
# -----------------------------------------------------------------
$t = $d.Tables.Item(1)
$nameCellPara = $t.Cell(1, 1).Range.Paragraphs.Item(1)
$nameXml = $nameCellPara.Range.WordOpenXML
$pPrStart = $nameXml.IndexOf("<w:pPr")
$pPrEnd = $nameXml.IndexOf("</w:pPr>", $pPrStart) + "</w:pPr>".Length
$namePPr = $nameXml.Substring($pPrStart, $pPrEnd - $pPrStart)

$newNameParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    $namePPr + `
    '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>NAME:</w:t></w:r>' + `
    '<w:r><w:rPr><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{e</w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>_</w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">name} </w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>{e</w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>_</w:t></w:r>' + `
    '<w:r><w:rPr><w:noProof/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>surname}</w:t></w:r>' + `
    '</w:p>'

$nameCellPara.Range.InsertXML($newNameParaXml)

# -----------------------------------------------------------------
# 3. ADDRESS / Email / Phone cells: "{{" -> "{" and "}}" -> "}".
#    Each " {{" / "}}" token lives entirely inside its own run in the
#    source document, so a plain text replace here leaves every run's
#    formatting/boundaries untouched - exactly matching the diff.
#    We rely on document order: after step 2 the NAME row no longer
#    contains any "{{"/"}}" tokens, so the first remaining " {{" (then
#    "}}") found from the top of the document is always ADDRESS, then
#    Email, then Phone, in that order.
# -----------------------------------------------------------------
Replace-FirstOccurrence " {{" " {"
Replace-FirstOccurrence " {{" " {"
Replace-FirstOccurrence " {{" " {"

Replace-FirstOccurrence "}}" "}"
Replace-FirstOccurrence "}}" "}"
Replace-FirstOccurrence "}}" "}"

# -----------------------------------------------------------------
# 4. Drop the "{% endfor %}" hidden runs.
# -----------------------------------------------------------------
Clear-HiddenControlParagraph "endfor"
